# Update workbook/sheet for new data date 2022-03-05 (commit: "Add data for 2022-03-13")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab name + workbook.xml <sheet name=.../>)
$ws.Name = "Through 2022-03-05"

# Update the "March (through 03-04)" label to "March (through 03-05)"
$ws.Range("A4").Value = "March (through 03-05)"

# Update March row (row 4) values for each year column B:I
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 11
$ws.Range("F4").Value = 6
$ws.Range("G4").Value = 9
$ws.Range("H4").Value = 13
$ws.Range("I4").Value = 30

# Update Total row (row 5) values for each year column B:I
$ws.Range("B5").Value = 41
$ws.Range("C5").Value = 93
$ws.Range("D5").Value = 136
$ws.Range("E5").Value = 148
$ws.Range("F5").Value = 85
$ws.Range("G5").Value = 150
$ws.Range("H5").Value = 355
$ws.Range("I5").Value = 331
